$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'codice_1_livello,codice_2_livello,label_ITA_1_livello,label_ENG_1_livello,label_ITA_2_livello,label_ENG_2_livello,references'
$ws.Cells.Item(2, 1).Value = '01,01.01,Disposizioni generali,General provisions,Piano triennale per la prevenzione della corruzione e della trasparenza,Plan for transparency and integrity,"D.lgs. n. 33/2013, Art. 10, c. 8, lett. a"'
$ws.Cells.Item(3, 1).Value = '01,01.02,Disposizioni generali,General provisions,Atti generali,General documents,"D.lgs. n. 33/2013, Art. 12, c. 1,2"'
$ws.Cells.Item(4, 1).Value = '01,01.03,Disposizioni generali,General provisions,Oneri informativi per cittadini e imprese,Obligatory information for citizens and companies,"D.lgs. n. 33/2013, Art. 34, c. 1,2,"'
$ws.Cells.Item(5, 1).Value = '02,02.01,Organizzazione,Organization,"Titolari di incarichi politici, di amministrazione, di direzione o di governo",Political and administrative bodies,"D.lgs. n. 33/2013, Art. 13, c. 1, lett. a | D.lgs. n. 33/2013, Art. 14"'
$ws.Cells.Item(6, 1).Value = '02,02.02,Organizzazione,Organization,Sanzioni per mancata comunicazione dei dati,Sanctions for failure of data communication,"D.lgs. n. 33/2013, Art. 47"'
$ws.Cells.Item(7, 1).Value = '02,02.03,Organizzazione,Organization,Rendiconti gruppi consiliari regionali/provinciali,Reports of regional / provincial council groups,"D.lgs. n. 33/2013, Art. 28, c. 1"'
$ws.Cells.Item(8, 1).Value = '02,02.04,Organizzazione,Organization,Articolazione degli uffici,Internal organization,"D.lgs. n. 33/2013, Art. 13, c. 1, lett. b,c"'
$ws.Cells.Item(9, 1).Value = '02,02.05,Organizzazione,Organization,Telefono o posta elettronica,Telephone and e-mail,"D.lgs. n. 33/2013, Art. 13, c. 1, lett. d"'
$ws.Cells.Item(10, 1).Value = '03,03.01,Consulenti e collaboratori,Consultants and collaborators,Titolari di incarichi di collaborazione o consulenza,Consultants and collaborators,"D.lgs. n. 33/2013, Art. 15, c. 1, lett. c"'
$ws.Cells.Item(11, 1).Value = '04,04.01,Personale,Personnel,Titolari di incarichi dirigenziali amministrativi di vertice,Top administrative tasks or offices,"D.lgs. n. 33/2013, Art. 15, c. 1,2 | D.lgs. n. 33/2013, Art. 41, c. 2,3"'
$ws.Cells.Item(12, 1).Value = '04,04.02,Personale,Personnel,Titolari di incarichi dirigenziali (dirigenti non generali),Executives,"D.lgs. n. 33/2013, Art. 14, c. 1, lett. a,b,c,d,e,f | D.lgs. n. 33/2013, Art. 14, c. 1ter | D.lgs. n. 33/2013, Art. 20, c. 3"'
$ws.Cells.Item(13, 1).Value = '04,04.03,Personale,Personnel,Dirigenti cessati,Executives ceased,"D.lgs. n. 33/2013, Art. 10, c. 8, lett. a,b,c,d,e,f"'
$ws.Cells.Item(14, 1).Value = '04,04.04,Personale,Personnel,Sanzioni per mancata comunicazione dei dati,Penalties for failure to communicate data ,"D.lgs. n. 33/2013, Art. 47, c. 1"'
$ws.Cells.Item(15, 1).Value = '04,04.05,Personale,Personnel,Posizioni organizzative,Organizational positions,"D.lgs. n. 33/2013, Art. 10, c. 8, lett. d"'
$ws.Cells.Item(16, 1).Value = '04,04.06,Personale,Personnel,Dotazione organica,Organic endowment,"D.lgs. n. 33/2013, Art. 16, c. 1,2"'
$ws.Cells.Item(17, 1).Value = '04,04.07,Personale,Personnel,Personale non a tempo indeterminato,Non-permanent staff,"D.lgs. n. 33/2013, Art. 17, c. 1,2"'
$ws.Cells.Item(18, 1).Value = '04,04.08,Personale,Personnel,Tassi di assenza,Absence rates,"D.lgs. n. 33/2013, Art. 16, c. 3"'
$ws.Cells.Item(19, 1).Value = '04,04.09,Personale,Personnel,Incarichi conferiti e autorizzati ai dipendenti (dirigenti e non dirigenti),Task assigned and authorized to employees (executives and non-executives),"D.lgs. n. 33/2013, Art. 18, c. 1"'
$ws.Cells.Item(20, 1).Value = '04,04.10,Personale,Personnel,Contrattazione collettiva,Collective negotiation,"D.lgs. n. 33/2013, Art. 21, c. 1"'
$ws.Cells.Item(21, 1).Value = '04,04.11,Personale,Personnel,Contrattazione integrativa,Integrative negotiation,"D.lgs. n. 33/2013, Art. 21, c. 2"'
$ws.Cells.Item(22, 1).Value = '04,04.12,Personale,Personnel,OIV (Organismo Indipendente di Valutazione),Independent Evaluation Agency,"D.lgs. n. 33/2013, Art. 10, c. 8, lett. c"'
$ws.Cells.Item(23, 1).Value = '05,05.01,Bandi di concorso,Competition announcements,Bandi di concorso,Competition announcements,"D.lgs. n. 33/2013, Art. 19"'
$ws.Cells.Item(24, 1).Value = '06,06.01,Performance,Performance,Sistema di misurazione e valutazione delle performance,Performance measurement and evaluation system,"Par. 1, delib. CIVIT n.104/2010"'
$ws.Cells.Item(25, 1).Value = '06,06.02,Performance,Performance,Piano delle Performance,Performance Plan,"D.lgs. n. 33/2013, Art. 10, c. 8, lett. b"'
$ws.Cells.Item(26, 1).Value = '06,06.03,Performance,Performance,Relazioni sulle Performance,Report on Performance,"D.lgs. n. 33/2013, Art. 10, c. 8, lett. b"'
$ws.Cells.Item(27, 1).Value = '06,06.04,Performance,Performance,Ammontare complessivo dei premi,Total amount of grants,"D.lgs. n. 33/2013, Art. 20, c. 1"'
$ws.Cells.Item(28, 1).Value = '06,06.05,Performance,Performance,Dati relativi ai premi,Informations on grants,"D.lgs. n. 33/2013, Art. 20, c. 2"'
$ws.Cells.Item(29, 1).Value = '07,07.01,Enti controllati,Controlled Institutions,Enti pubblici vigilati,Supervised public units,"D.lgs. n. 33/2013, Art. 22, c. 1, lett. a | D.lgs. n. 33/2013, Art. 22, c. 2,3"'
$ws.Cells.Item(30, 1).Value = '07,07.02,Enti controllati,Controlled Institutions,SocietÃ  partecipate,Subsidiaries companies,"D.lgs. n. 33/2013, Art. 22, c. 1, lett. b | D.lgs. n. 33/2013, Art. 22, c. 2,3"'
$ws.Cells.Item(31, 1).Value = '07,07.03,Enti controllati,Controlled Institutions,Enti di diritto privato controllati,Controlled private law units,"D.lgs. n. 33/2013, Art. 22, c. 1, lett. c | D.lgs. n. 33/2013, Art. 22, c. 2,3"'
$ws.Cells.Item(32, 1).Value = '07,07.04,Enti controllati,Controlled Institutions,Rappresentazione grafica,Organisation chart,"D.lgs. n. 33/2013, Art. 22, c. 1, lett. d"'
$ws.Cells.Item(33, 1).Value = '08,08.01,AttivitÃ  e procedimenti,Activity and procedures,Tipologie di procedimento,Procedures types,"D.lgs. n. 33/2013, Art. 35, c. 1,2"'
$ws.Cells.Item(34, 1).Value = '08,08.02,AttivitÃ  e procedimenti,Activity and procedures,Dichiarazioni sostitutive e acquisizione d''ufficio dei dati,Substitutive declarations and data acquisition ex officio,"D.lgs. n. 33/2013, Art. 35, c. 3"'
$ws.Cells.Item(35, 1).Value = '09,09.01,Provvedimenti,Measures,Provvedimenti organi indirizzo politico,Provisions adopted by the political direction bodies,"D.lgs. n. 33/2013, Art. 23, c. 1"'
$ws.Cells.Item(36, 1).Value = '09,09.02,Provvedimenti,Measures,Provvedimenti dirigenti amministrativi,Provisions adopted by administrative managers,"D.lgs. n. 33/2013, Art. 23, c. 1"'
$ws.Cells.Item(37, 1).Value = '10,10.01,Bandi di gara e contratti,Competition announcements and contracts,Informazioni sulle singole procedure in formato tabellare,Single procedure information in tabular format,"Art. 4, delib. ANAC n.39/2016 | D.lgs. n. 33/2013, Art. 1, c. 32 | D.lgs. n. 33/2013, Art. 37, c. 1, lett. a"'
$ws.Cells.Item(38, 1).Value = '10,10.02,Bandi di gara e contratti,Competition announcements and contracts,Atti delle amministrazioni aggiudicatrici e degli enti aggiudicatori distintamente per ogni procedura,Acts of contracting authorities and contracting entities separately for each procedure,"D.lgs. n. 33/2013, Art. 37, c. 1, lett. b"'
$ws.Cells.Item(39, 1).Value = '11,11.01,"Sovvenzioni, contributi, sussidi, vantaggi economici","Subsidies, contributions, grants, economic advantages",Criteri e modalitÃ ,Criteria and methods,"D.lgs. n. 33/2013, Art. 26, c. 1"'
$ws.Cells.Item(40, 1).Value = '11,11.02,"Sovvenzioni, contributi, sussidi, vantaggi economici","Subsidies, contributions, grants, economic advantages",Atti di concessione,Concession acts,"D.lgs. n. 33/2013, Art. 26, c. 2 | D.lgs. n. 33/2013, Art. 27"'
$ws.Cells.Item(41, 1).Value = '12,12.01,Bilanci,Budgets,Bilancio preventivo e consuntivo,Budget and final balance sheet,"D.lgs. n. 33/2013, Art. 29, c. 1"'
$ws.Cells.Item(42, 1).Value = '12,12.02,Bilanci,Budgets,Piano degli indicatori e dei risultati attesi di bilancio,Plan of indicators and expected balance sheet results,"D.lgs. n. 33/2013, Art. 29, c. 2"'
$ws.Cells.Item(43, 1).Value = '13,13.01,Beni immobili e gestione patrimonio,Real estate and capital treatment,Patrimonio immobiliare,Real estate,"D.lgs. n. 33/2013, Art. 30"'
$ws.Cells.Item(44, 1).Value = '13,13.02,Beni immobili e gestione patrimonio,Real estate and capital treatment,Canoni di locazione o affitto,Rents or renting,"D.lgs. n. 33/2013, Art. 30"'
$ws.Cells.Item(45, 1).Value = '14,14.01,Controlli e rilievi sull''amministrazione,Checks and measurements about administration,"Organismi dipendenti di valutazione, nuclei di valutazione o altri organismi con funzioni analoghe","Evaluation dependent bodies, evaluation units or other bodies with similar functions","D.lgs. n. 33/2013, Art. 31"'
$ws.Cells.Item(46, 1).Value = '14,14.02,Controlli e rilievi sull''amministrazione,Checks and measurements about administration,Organi di revisione amministrativa e contabile,Administrative and accounting review bodies,"D.lgs. n. 33/2013, Art. 31"'
$ws.Cells.Item(47, 1).Value = '14,14.03,Controlli e rilievi sull''amministrazione,Checks and measurements about administration,Corte dei conti,Corte dei conti,"D.lgs. n. 33/2013, Art. 31"'
$ws.Cells.Item(48, 1).Value = '15,15.01,Servizi erogati,Services provided,Carta dei servizi e standard di qualitÃ ,Chart of services and quality standards,"D.lgs. n. 33/2013, Art. 32, c. 1"'
$ws.Cells.Item(49, 1).Value = '15,15.02,Servizi erogati,Services provided,Class action,Class action,"D.lgs. n. 198/2009, Art. 1, c. 2 | D.lgs. n. 198/2009, Art. 4, c. 2,6"'
$ws.Cells.Item(50, 1).Value = '15,15.03,Servizi erogati,Services provided,Costi contabilizzati,Recorded costs,"D.lgs. n. 33/2013, Art. 10, c. 5 | D.lgs. n. 33/2013, Art. 32, c. 2, lett. a"'
$ws.Cells.Item(51, 1).Value = '15,15.04,Servizi erogati,Services provided,Liste di attesa,Waiting lists,"D.lgs. n. 33/2013, Art. 41, c. 6,"'
$ws.Cells.Item(52, 1).Value = '15,15.05,Servizi erogati,Services provided,Servizi in rete,Network services,"D.lgs. n. 179/2016, Art. 8, c. 1"'
$ws.Cells.Item(53, 1).Value = '15,15.06,Servizi erogati,Services provided,Dati sui pagamenti,Payments data,"D.lgs. n. 33/2013, Art. 4bis, c. 2"'
$ws.Cells.Item(54, 1).Value = '15,15.07,Servizi erogati,Services provided,Dati sui pagamenti del servizio sanitario nazionale,Data on payments from the national health service,"D.lgs. n. 33/2013, Art. 36"'
$ws.Cells.Item(55, 1).Value = '15,15.08,Servizi erogati,Services provided,Indicatore di tempestivitÃ  dei pagamenti,Average times of service provision,"D.lgs. n. 33/2013, Art. 33"'
$ws.Cells.Item(56, 1).Value = '15,15.09,Servizi erogati,Services provided,IBAN e pagamenti informatici,IBAN and electronic payments,"D.lgs. n. 33/2013, Art. 36"'
$ws.Cells.Item(57, 1).Value = '16,16.01,Opere pubbliche,Public works,Nuclei di valutazione eÂ verifica degli investimenti pubblici,Evaluation and verification of public investments,"D.lgs. n. 33/2013, Art. 38, c. 1"'
$ws.Cells.Item(58, 1).Value = '16,16.02,Opere pubbliche,Public works,Atti di programmazione delle opere pubbliche,Planning of public works,"D.lgs. n. 33/2013, Art. 38, c. 2"'
$ws.Cells.Item(59, 1).Value = '16,16.03,Opere pubbliche,Public works,Tempi costi e indicatori di realizzazione delle opere pubbliche,"Cost, times and indicators for the realization of public works","D.lgs. n. 33/2013, Art. 38, c. 2"'
$ws.Cells.Item(60, 1).Value = '17,17.01,Pianificazione e governo del territorio,Planning and government of the area,Pianificazione e governo del territorio,Planning and government of the area,"D.lgs. n. 33/2013, Art. 39"'
$ws.Cells.Item(61, 1).Value = '18,18.01,Informazioni ambientali,Environmental information,Informazioni ambientali,Environmental information,"D.lgs. n. 33/2013, Art. 40"'
$ws.Cells.Item(62, 1).Value = '19,19.01,Strutture sanitarie private accreditate,Health accredited facilities,Strutture sanitarie private accreditate,Health accredited facilities,"D.lgs. n. 33/2013, Art. 41, c. 4"'
$ws.Cells.Item(63, 1).Value = '20,20.01,Interventi straordinari e di emergenza,Extraordinary and emergency actions,Interventi straordinari e di emergenza,Extraordinary and emergency actions,"D.lgs. n. 33/2013, Art. 42"'
$ws.Cells.Item(64, 1).Value = '21,21.01,Altri contenuti,Other contents,Prevenzione della Corruzione,Prevention of corruption,"D.lgs. n. 33/2013, Art. 1, c. 8 | D.lgs. n. 33/2013, Art. 10, c. 8, lett. a | D.lgs. n. 33/2013, Art. 18, c. 5"'
$ws.Cells.Item(65, 1).Value = '21,21.02,Altri contenuti,Other contents,Accesso civico,Civic access,"D.lgs. n. 33/2013, Art. 5, c. 1,2"'
$ws.Cells.Item(66, 1).Value = '21,21.03,Altri contenuti,Other contents,"AccessibilitÃ  e Catalogo dei dati, metadati e banche dati","Accessibility and data catalog, metadata and databases","D.lgs. n. 33/2013, Art. 53, c. 1bis | D.lgs. n. 33/2013, Art. 9, c. 7"'
$ws.Cells.Item(67, 1).Value = '21,21.04,Altri contenuti,Other contents,Dati ulteriori,Other informations,"D.lgs. n. 33/2013, Art. 7bis, c. 3"'

$ws.Range("A68:A70").ClearContents()

